# Update issue 260 update performance document.
# Adds a new "v1433" performance sample column to both the "Sponza"
# worksheet (new column M) and the "ComplexMesh" worksheet (new column L),
# mirroring the layout/formulas of the existing rightmost data column on
# each sheet.

$wb = $excel.ActiveWorkbook
$wsSponza  = $wb.Worksheets.Item("Sponza")
$wsComplex = $wb.Worksheets.Item("ComplexMesh")

# ---------------------------------------------------------------------
# Sheet "Sponza": new column M (copy layout/style from column L)
# ---------------------------------------------------------------------

$wsSponza.Range("L1:L16").Copy($wsSponza.Range("M1"))

$wsSponza.Range("M1").Value = "v1433"

$sponzaValues = @(7019,7019,6970,6983,6982,6956,6984,6982,7001,6972)
for ($i = 0; $i -lt $sponzaValues.Length; $i++) {
  $wsSponza.Cells.Item(2 + $i, 13).Value = $sponzaValues[$i]
}

$wsSponza.Range("M12").Formula = "=AVERAGE(M2:M11)"
$wsSponza.Range("M13").Formula = "=_xlfn.VAR.S(M2:M11)"
$wsSponza.Range("M14").Formula = "=1-_xlfn.T.TEST(L2:L11,M2:M11,2,3)"
$wsSponza.Range("M15").Formula = "=L12/M12"
$wsSponza.Range("M16").Formula = "=B12/M12"

# Extend the conditional formatting range B15:L16 -> B15:M16
$sponzaNewRange = $wsSponza.Range("B15:M16")
$sponzaRules = $wsSponza.Range("B15:L16").FormatConditions
for ($i = 1; $i -le $sponzaRules.Count; $i++) {
  $sponzaRules.Item($i).ModifyAppliesToRange($sponzaNewRange)
}

# ---------------------------------------------------------------------
# Sheet "ComplexMesh": new column L (copy layout/style from column K)
# ---------------------------------------------------------------------

$wsComplex.Range("K1:K16").Copy($wsComplex.Range("L1"))

$wsComplex.Range("L1").Value = "v1433"

$complexValues = @(4917,4967,5046,4977,4964,4943,4952,4947,4949,4944)
for ($i = 0; $i -lt $complexValues.Length; $i++) {
  $wsComplex.Cells.Item(2 + $i, 12).Value = $complexValues[$i]
}

$wsComplex.Range("L12").Value = 4964
$wsComplex.Range("L13").Formula = "=_xlfn.VAR.S(L2:L11)"
$wsComplex.Range("L14").Formula = "=1-_xlfn.T.TEST(K2:K11,L2:L11,2,3)"
$wsComplex.Range("L15").Formula = "=K12/L12"
$wsComplex.Range("L16").Formula = "=B12/L12"

# Extend the conditional formatting range B15:K16 -> B15:L16
$complexNewRange = $wsComplex.Range("B15:L16")
$complexRules = $wsComplex.Range("B15:K16").FormatConditions
for ($i = 1; $i -le $complexRules.Count; $i++) {
  $complexRules.Item($i).ModifyAppliesToRange($complexNewRange)
}

# ---------------------------------------------------------------------
# Selections: update the active cell on each sheet, and restore
# "ComplexMesh" as the selected/active tab (matching the original file).
# ---------------------------------------------------------------------

$wsSponza.Activate()
$wsSponza.Range("M2").Select()

$wsComplex.Activate()
$wsComplex.Range("L2").Select()
